# Update cryptos list values to reflect latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells hold text-formatted numbers (e.g. "233.00", "1.814.06").
# Force the number format to Text before assigning so Excel does not silently
# reinterpret/round them as floating point numbers and lose the exact digits.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.982.23"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.814.06"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").Value = "233.00"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("E6").Value = "  +1.82%  "
$ws.Range("D7").Value = "1.01"
$ws.Range("D8").Value = "40.20"
$ws.Range("E8").Value = "  -10.95%  "
$ws.Range("D9").Value = "0.321"
$ws.Range("E9").Value = "  +8.45%  "
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "2.077.07"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "1.810.16"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "11.10"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("D17").Value = "34.976.40"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "69.42"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "238.12"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").Value = "4.65"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("D25").Value = "172.72"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "7.82"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "17.42"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  +32.18%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").Value = "3.340.71"
$ws.Range("E31").Value = "  +37.50%  "
$ws.Range("D32").Value = "0.0554"
$ws.Range("E32").Value = "  +6.68%  "
$ws.Range("D33").Value = "3.93"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "3.96"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("B36").Value = "Aave"
$ws.Range("C36").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D36").Value = "93.03"
$ws.Range("E36").Value = "  +4.60%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.13"
$ws.Range("E37").Value = "  +7.01%  "
$ws.Range("E38").Value = "  +3.47%  "
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").Value = "1.28"
$ws.Range("E40").Value = "  +5.04%  "
$ws.Range("D41").Value = "1.306.27"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("D42").Value = "0.984"
$ws.Range("E42").Value = "  +2.91%  "
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "14.63"
$ws.Range("E44").Value = "  -4.17%  "
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("D47").Value = "6.29"
$ws.Range("E47").Value = "  +6.24%  "
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "1.992.86"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "0.0646"
$ws.Range("E51").Value = "  +5.78%  "
